$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new row 8 for P1-006 (identifier / date / firmware version) ---
$ws.Range("A8").Value = "P1-006"

$ws.Range("B8").Value = 44520
$ws.Range("B8").NumberFormat = "d-mmm"

$ws.Range("C8").Value = "pre-1.0.6"

# --- Update row 3 (P1-001): bump firmware version and append a new history line ---
$ws.Range("C3").Value = "pre-1.0.6"

$existingNote = $ws.Range("E3").Value2
$ws.Range("E3").Value = $existingNote + "`n1/23: Upgraded firmware to pre-1.0.6"
$ws.Range("E3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 34

# --- Finish row 8 with the new history note ---
$ws.Range("E8").Value = "1/23: Upgraded firmware to pre-1.0.6"

$ws.Range("E8").Select() | Out-Null
